$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B86 was an inline string "2"; change it to a real number 2
$ws.Range("B86").Value = 2

# Append a new row 87 with the new annotation data for Ruilin
$ws.Range("A87").Value = "Ruilin"

# B87 must stay a text value "3" (not auto-converted to a number)
$b87 = $ws.Range("B87")
$b87.NumberFormat = "@"
$b87.Value = "3"
$b87.ClearFormats()

$ws.Range("C87").Value = "无"
$ws.Range("D87").Value = "DFT"
$ws.Range("E87").Value = "MET"
$ws.Range("F87").Value = "94664fc5-740b-497e-9f27-9fbb6b5fbbdd"
$ws.Range("G87").Value = "TT0bFo9VZpFWg_annotated.xlsx"
$ws.Range("H87").Value = "The net gets bigger, yet keeps underfitting the training set."
